$d = $word.ActiveDocument

$oldTitle = "Balanza de pagos: Explorando los sistemas monetarios internacionales, tipos de cambio fijos y flexibles, y el papel en evolución del FMI en la economía mundial."
$newTitle = "Estructura de la balanza de pagos"

# Replace both occurrences of the Heading1 title text (the "title" bookmark
# heading near the top, and the duplicate "firstheader" bookmark heading
# further down) by rewriting each matching paragraph's text in place.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ($oldTitle + "`r")) {
        $p.Range.Text = $newTitle
    }
}

# Locate the "Palabras clave" paragraph inside the Resumen/abstract section
# and its index, then insert a new placeholder abstract paragraph right
# before it (i.e. right after the "Resumen" Heading1 paragraph).
$kwIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Palabras clave*") {
        $kwIndex = $i
        break
    }
}

$resumenPara = $d.Paragraphs.Item($kwIndex - 1)

$insertPoint = $d.Range($resumenPara.Range.End - 1, $resumenPara.Range.End - 1)
$insertPoint.InsertAfter("`rEste abstract será actualizado una vez que se complete el contenido final del artículo.")

$newPara = $d.Paragraphs.Item($kwIndex)
$newPara.Style = "AbstractFirstParagraph"

$kwPara = $d.Paragraphs.Item($kwIndex + 1)
$kwPara.Style = "BodyText"
